$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the remaining fields of the existing "Move Cases" row (row 61).
$ws.Range("B61").Value = "move 28 cases to Closed/Archived"
$ws.Range("D61").Value = 44119
$ws.Range("E61").Value = "GSTC"
$ws.Range("F61").Value = "Done"

# Add a new "Move Cases" row (row 62).
$ws.Range("A62").Value = "Move Cases"
$ws.Range("B62").Value = "move 51 cases to Closed/Archived"
$ws.Range("D62").Value = 44119
$ws.Range("E62").Value = "GSTC"
$ws.Range("F62").Value = "Done"

# Add a new "reforward correspondence" row (row 63).
$ws.Range("A63").Value = "reforward correspondence "
$ws.Range("B63").Value = "reforword correspondence "
$ws.Range("D63").Value = 44119
$ws.Range("E63").Value = "Murasalat"
$ws.Range("F63").Value = "Done"

# Copy the formatting (only) of row 61 down onto the two newly added rows,
# and match its row height, then resize the Table42 list object /
# autofilter to cover the new rows.
$ws.Range("A61:F61").Copy()
$ws.Range("A62:F63").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Rows.Item(62).RowHeight = 24.95
$ws.Rows.Item(63).RowHeight = 24.95

$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:F63"))

# Update the view to match the edited workbook state.
$ws.Range("F64").Select()
